$d = $word.ActiveDocument

# wdReplace constant used throughout: 2 = wdReplaceAll
# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)

# --- Whole-phrase replacements (each phrase is a single cached
#     MERGEFIELD result run, occurring exactly twice in the body) ---

$d.Content.Find.Execute("MR", $true, $true, $false, $false, $false, `
    $true, 1, $false, "MS", 2)

$d.Content.Find.Execute("GARRY L. BAYOT", $true, $false, $false, $false, `
    $false, $true, 1, $false, "FELICITAS M. SUMAGUI", 2)

$d.Content.Find.Execute("Vice Mayor'S Office", $true, $false, $false, `
    $false, $false, $true, 1, $false, "City Social Welfare Development Office", 2)

$d.Content.Find.Execute("September 01, 2003", $true, $false, $false, `
    $false, $false, $true, 1, $false, "November 10, 2003", 2)

$d.Content.Find.Execute("his resignation", $true, $false, $false, $false, `
    $false, $true, 1, $false, "her Compulsory retirement", 2)

$d.Content.Find.Execute("July 01, 2022", $true, $false, $false, $false, `
    $false, $true, 1, $false, "March 23, 2023", 2)

$d.Content.Find.Execute("June", $true, $true, $false, $false, $false, `
    $true, 1, $false, "July", 2)

# --- Narrow, position-scoped replacements for the single-character day
#     number ("2" -> "5") and its ordinal suffix ("nd" -> "th"). These are
#     too short/generic for a safe whole-document Find, so each match is
#     located via a bounded Range anchored right after the unique phrase
#     "issued this " (there are two such occurrences, one per certificate
#     copy in the body), which keeps the surrounding run formatting
#     (including the superscript on the suffix) untouched.

$searchStart = 0
for ($i = 0; $i -lt 2; $i++) {
    $anchor = $d.Range($searchStart, $d.Content.End)
    $anchor.Find.Execute("issued this ")
    $anchor.Collapse(0)
    $anchor.MoveEnd(1, 40)

    $dayRange = $d.Range($anchor.Start, $anchor.End)
    $dayRange.Find.Execute("2")
    $dayRange.Text = "5"

    $suffixRange = $d.Range($anchor.Start, $anchor.End)
    $suffixRange.Find.Execute("nd")
    $suffixRange.Text = "th"

    $searchStart = $anchor.End
}
